$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update G5 value 262 -> 1262
$ws.Range("G5").Value = 1262

# Copy formatting from row 7 to row 8 so new row matches existing style pattern
$ws.Range("A7:L7").Copy()
$ws.Range("A8:L8").PasteSpecial(-4122)

# Fill in row 8 data
$ws.Range("B8").Value = 24
$ws.Range("C8").Formula = "=3*38835"
$ws.Range("D8").Formula = "=C8*F8"
$ws.Range("E8").Value = 277
$ws.Range("F8").Value = 1.5
$ws.Range("G8").Value = 1411
$ws.Range("H8").Value = 0.042418981481481481
$ws.Range("I8").Value = 6700
$ws.Range("J8").Value = "Vampiro"
$ws.Range("K8").Value = "Normal"
$ws.Range("L8").Value = 46013

# Update selection to G6
$ws.Range("G6").Select()
